# Evaluation-data update: add new detections for rows 5 & 6, revise row 2.
# (see commit message: "luminosity" filtering improved the accuracy figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: revised detection counts / success rate ---
$ws.Range("B2").Value = 0.1
$ws.Range("D2").Value = 10
$ws.Range("F2").Value = 9

# --- Row 5: new matching-points data (previously blank but for A5/G5/H5) ---
$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 5.5
$ws.Range("F5").Value = 7

# --- Row 6: new matching-points data (previously blank but for A6/G6/H6) ---
$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 6

# Leave the selection where the author left off editing.
[void]$ws.Range("C6").Select()
